$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the new date as literal text (not an auto-converted date serial) by
# routing it through a scratch formula cell and pasting its computed value.
# This preserves the existing shared-string / number-format styling used by
# the rest of column A instead of minting a brand new cell style.
$ws.Range("Z1").Formula = "=""2021/10/04"""
$ws.Range("Z1").Copy()
$ws.Range("A89").PasteSpecial(-4163)
$ws.Range("Z1").Clear()

$ws.Range("B89").Value = 94.3
$ws.Range("C89").Value = 95.1
$ws.Range("D89").Value = 0.91
$ws.Range("E89").Value = 0.9

$ws.Range("A90").Select()
